$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 13 ("Programa resumido:" and everything below
# shifts down by one). This makes room for a dedicated "Docentes
# responsáveis:" answer row, mirroring the two preceding rows
# (Objetivos:/Objectives:) that already carry their own answer rows.
$ws.Rows.Item(13).Insert()

# Row 10 "Objetivos:" previously (incorrectly) duplicated the teacher name
# in B/C; it should hold the Portuguese objectives paragraph instead.
$ws.Range("B10").Value = "Complementar a formação dos estudantes abordando, com maior profundidade, tópicos atuais e relevantes e atualizar com temas no estado da arte."
$ws.Range("C10").Value = "Complementar a formação dos estudantes abordando, com maior profundidade, tópicos atuais e relevantes e atualizar com temas no estado da arte."

# New row 13: the "Docentes responsáveis:" answer (no label in column A, so
# drop the blank carried-over A13 cell). B/C reuse the answer-cell look from
# row 11 (normal wrap text / red wrap text) instead of the bold label style
# the row-insert copied down from column A.
$ws.Range("A13").Clear()
$ws.Range("B11").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C11").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Range("B13").Value = "11079086 - Herlandí de Souza Andrade"
$ws.Range("C13").Value = "11079086 - Herlandí de Souza Andrade"

# Row 14 ("Programa resumido:") previously held "Semestral"; update to the
# correct short-syllabus answer.
$ws.Range("B14").Value = "A definir de acordo com o tópico programado"
$ws.Range("C14").Value = "A definir de acordo com o tópico programado"

# Row 16 ("Programa:") previously (incorrectly) held the activation date;
# fill in the real syllabus paragraph.
$ws.Range("B16").Value = "O conteúdo desta disciplina optativa será de acordo com o tópico a ser programado, devendo abordar assuntos complementares ao conteúdo regular do curso de graduação."
$ws.Range("C16").Value = "O conteúdo desta disciplina optativa será de acordo com o tópico a ser programado, devendo abordar assuntos complementares ao conteúdo regular do curso de graduação."

# Row 19 ("Método:") previously (incorrectly) duplicated the teacher name;
# fill in the real evaluation-method paragraph.
$ws.Range("B19").Value = "Esta disciplina deverá conter no mínimo duas avaliações denominadas A1 e A2. As avalições poderão ser: escritas, práticas, seminários, trabalhos de campo, projetos, ou outra forma de avaliação definida pelo professor."
$ws.Range("C19").Value = "Esta disciplina deverá conter no mínimo duas avaliações denominadas A1 e A2. As avalições poderão ser: escritas, práticas, seminários, trabalhos de campo, projetos, ou outra forma de avaliação definida pelo professor."

# Row 20 ("Critério:") gets the weighted-average criterion text.
$ws.Range("B20").Value = "Média ponderada das avaliações (M)."
$ws.Range("C20").Value = "Média ponderada das avaliações (M)."

# Row 21 ("Norma de recuperação:") gets the make-up-exam rule text.
$ws.Range("B21").Value = "A recuperação será composta por uma única prova (RC) englobando toda a matéria ministrada ao longo do semestre. A média final, para os alunos em recuperação, será calculada com base na relação: MF=(M+RC)/2"
$ws.Range("C21").Value = "A recuperação será composta por uma única prova (RC) englobando toda a matéria ministrada ao longo do semestre. A média final, para os alunos em recuperação, será calculada com base na relação: MF=(M+RC)/2"

# Row 22 ("Bibliografia:") previously held the make-up-exam text; replace
# with the actual bibliography paragraph.
$ws.Range("B22").Value = "Livros, artigos ou texto fornecido pelo docente responsável extraídos de livros ou revistas especializadas na área de Engenharia de Produção."
$ws.Range("C22").Value = "Livros, artigos ou texto fornecido pelo docente responsável extraídos de livros ou revistas especializadas na área de Engenharia de Produção."

# Column layout: column A is now only 1 column wide (was merged 1:2 before).
$ws.Columns.Item(1).ColumnWidth = 30.7109375
